# Add new column R ("l1") to Sheet1, mirroring the P/Q "c1" style columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold + border, style index 1) from Q1 onto R1
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)

# Set the new header text and the zero values for every data row (2-79)
$ws.Range("R1").Value = "l1"
$ws.Range("R2:R79").Value = 0

Write-Host "done"
